$d = $word.ActiveDocument

function Replace-Span($oldText, $newText) {
    $full = $d.Content.Text
    $pos = $full.IndexOf($oldText)
    if ($pos -lt 0) {
        throw "Could not find text: $oldText"
    }
    $r = $d.Range($pos, $pos + $oldText.Length)
    $r.Text = $newText
}

# ---------------------------------------------------------------------------
# 1) Paragraph 2 ("Life is all about healthy relations...") - remove the
#    spell-check proofErr wrappers around "KhalsaAid"/"SikhNet" by forcing a
#    run merge across that stretch of text (text itself is unchanged).
# ---------------------------------------------------------------------------
$old1 = "the fundraisers initiated by the KhalsaAid and SikhNet that have organized relief camps and rescues for the victims of wars and natural calamities. These fundraisers can be utilized for the local growth and development. Many a time, the local news has never been able to reach to media. In such a scenario, Facebook like platforms can spread the word among people and keep them up to date with the information."
Replace-Span $old1 $old1

# ---------------------------------------------------------------------------
# 2) Paragraph 7 ("The Facebook Local app...") - "convenient" -> "comfortable"
#    and "the messenger" -> "the Facebook messenger".
# ---------------------------------------------------------------------------
$old2 = "everyone feels convenient to be accessible via phone. However, the messenger acts as a savior and gives a wider reach."
$new2 = "everyone feels comfortable to be accessible via phone. However, the Facebook messenger acts as a savior and gives a wider reach."
Replace-Span $old2 $new2

# ---------------------------------------------------------------------------
# 3) Paragraph 10 - several word-level edits:
#    "housewife" -> "homemaker"
#    "used to enthrall" -> "enthralled"
#    "which lead me" -> "which led me"
#    "I was the first one" -> "I am the first one"
# ---------------------------------------------------------------------------
$old3 = "Ostensibly, my mother is a housewife. My parents always wanted me to lead a facile life, completing my Bachelor's in Arts and eventually get married. However, I was always enamored of programming and technical stuff and all the innovations of Google, Microsoft, and other corporations used to enthrall me. From the time of my high school, I dreamt of being a software engineer and eventually, have my software company and be an independent businesswoman which lead me to make antithetical choices. Though I had to compromise with a Tier 3 engineering school, I was the first one in the family to pursue a degree in Computer Science."
$new3 = "Ostensibly, my mother is a homemaker. My parents always wanted me to lead a facile life, completing my Bachelor's in Arts and eventually get married. However, I was always enamored of programming and technical stuff and all the innovations of Google, Microsoft, and other corporations enthralled me. From the time of my high school, I dreamt of being a software engineer and eventually, have my software company and be an independent businesswoman which led me to make antithetical choices. Though I had to compromise with a Tier 3 engineering school, I am the first one in the family to pursue a degree in Computer Science."
Replace-Span $old3 $new3

# ---------------------------------------------------------------------------
# 4) Paragraph 12 - rewrite of the "woman this seemed..." sentence, plus
#    appended new closing sentence, and removal of the Systers proofErr tags.
# ---------------------------------------------------------------------------
$old4 = "Being a woman this seemed to be another obstacle in the sense that women might not be accepted in the industry. However, I had to take this challenge and transpose the trend. This also influenced me to be a part of Systers, an Anita Borg Institute community. Systers welcomes women in technology to collaborate and share their experiences. Furthermore, I was one of the recipients of Google Scholarship for 2017 Grace Hopper Celebration."
$new4 = "Being a woman this seems to be another obstacle in the sense that women are not accepted in this industry. However, I have taken this challenge to transpose the trend. This also influenced me to be a part of Systers, an Anita Borg Institute community. Systers welcomes women in technology to collaborate and share their experiences. Furthermore, I was one of the recipients of Google Scholarship for 2017 Grace Hopper Celebration. All this motivates me to continue to pave way towards my goals and reminding myself that these are attainable."
Replace-Span $old4 $new4

# ---------------------------------------------------------------------------
# 5) Move the "_GoBack" bookmark from before "Not to forget..." (paragraph 7)
#    to its new location near the end of the document, between "...reminding
#    myself that these are" and " attainable." (paragraph 12).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$fullText = $d.Content.Text
$anchor = "reminding myself that these are"
$pos = $fullText.IndexOf($anchor) + $anchor.Length
$bookmarkRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
